{"js": "const body = context.document.body;\n\nconst pairs = [\n  [\"2025-01-08 Wednesday\", \"2025-01-14 Tuesday\"],\n  [\"78\u00d751=3978\", \"39\u00d785=3315\"],\n  [\"59\u00d711=649\", \"92\u00d774=6808\"],\n  [\"22\u00d748=1056\", \"97\u00d747=4559\"],\n  [\"77\u00d762=4774\", \"78\u00d721=1638\"],\n  [\"78\u00d738=2964\", \"84\u00d730=2520\"],\n  [\"37\u00d748=1776\", \"95\u00d797=9215\"],\n  [\"85\u00d750=4250\", \"57\u00d719=1083\"],\n  [\"79\u00d771=5609\", \"71\u00d755=3905\"],\n  [\"16\u00d786=1376\", \"62\u00d771=4402\"],\n  [\"81\u00d712=972\", \"99\u00d754=5346\"],\n  [\"69\u00d773=5037\", \"14\u00d750=700\"],\n  [\"96\u00d798=9408\", \"41\u00d747=1927\"],\n  [\"64\u00d795=6080\", \"86\u00d776=6536\"],\n  [\"89\u00d744=3916\", \"82\u00d711=902\"],\n  [\"36\u00d783=2988\", \"56\u00d763=3528\"],\n  [\"47\u00d742=1974\", \"45\u00d756=2520\"],\n  [\"29\u00d736=1044\", \"23\u00d713=299\"],\n  [\"21\u00d731=651\", \"99\u00d788=8712\"],\n  [\"22\u00d765=1430\", \"93\u00d754=5022\"],\n  [\"22\u00d788=1936\", \"31\u00d733=1023\"],\n  [\"12\u00d787=1044\", \"11\u00d764=704\"],\n  [\"31\u00d785=2635\", \"35\u00d793=3255\"],\n  [\"73\u00d756=4088\", \"56\u00d734=1904\"],\n  [\"46\u00d749=2254\", \"45\u00d787=3915\"],\n  [\"44\u00d789=3916\", \"51\u00d748=2448\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-01-08 Wednesday\", \"2025-01-14 Tuesday\"),\n    @(\"78\u00d751=3978\", \"39\u00d785=3315\"),\n    @(\"59\u00d711=649\", \"92\u00d774=6808\"),\n    @(\"22\u00d748=1056\", \"97\u00d747=4559\"),\n    @(\"77\u00d762=4774\", \"78\u00d721=1638\"),\n    @(\"78\u00d738=2964\", \"84\u00d730=2520\"),\n    @(\"37\u00d748=1776\", \"95\u00d797=9215\"),\n    @(\"85\u00d750=4250\", \"57\u00d719=1083\"),\n    @(\"79\u00d771=5609\", \"71\u00d755=3905\"),\n    @(\"16\u00d786=1376\", \"62\u00d771=4402\"),\n    @(\"81\u00d712=972\", \"99\u00d754=5346\"),\n    @(\"69\u00d773=5037\", \"14\u00d750=700\"),\n    @(\"96\u00d798=9408\", \"41\u00d747=1927\"),\n    @(\"64\u00d795=6080\", \"86\u00d776=6536\"),\n    @(\"89\u00d744=3916\", \"82\u00d711=902\"),\n    @(\"36\u00d783=2988\", \"56\u00d763=3528\"),\n    @(\"47\u00d742=1974\", \"45\u00d756=2520\"),\n    @(\"29\u00d736=1044\", \"23\u00d713=299\"),\n    @(\"21\u00d731=651\", \"99\u00d788=8712\"),\n    @(\"22\u00d765=1430\", \"93\u00d754=5022\"),\n    @(\"22\u00d788=1936\", \"31\u00d733=1023\"),\n    @(\"12\u00d787=1044\", \"11\u00d764=704\"),\n    @(\"31\u00d785=2635\", \"35\u00d793=3255\"),\n    @(\"73\u00d756=4088\", \"56\u00d734=1904\"),\n    @(\"46\u00d749=2254\", \"45\u00d787=3915\"),\n    @(\"44\u00d789=3916\", \"51\u00d748=2448\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
